$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("deliveries")

# Add a new delivery record in row 5
$ws.Cells.Item(5, 1).Value = 4

# Copy the date format from the row above so the new date cell keeps the
# same "m/d/yyyy" style (s="1") instead of creating a new numFmt.
$ws.Cells.Item(4, 2).Copy() | Out-Null
$ws.Cells.Item(5, 2).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(5, 2).Value = "1/10/2025"

$ws.Cells.Item(5, 4).Value = 1555
$ws.Cells.Item(5, 6).Value = "pista"
$ws.Cells.Item(5, 7).Value = "delihvery"
$ws.Cells.Item(5, 5).Value = "15kg"

# Update the selected cell to match the saved view state
$ws.Range("F14").Select() | Out-Null
